$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "67.924.41"
Set-TextValue 2 5 "  -0.63%  "
Set-TextValue 3 4 "3.812.11"
Set-TextValue 3 5 "  -1.87%  "
Set-TextValue 4 4 "0.998"
Set-TextValue 4 5 "  -0.04%  "
Set-TextValue 5 4 "599.94"
Set-TextValue 5 5 "  -0.37%  "
Set-TextValue 6 4 "168.95"
Set-TextValue 6 5 "  +0.53%  "
Set-TextValue 7 4 "3.812.45"
Set-TextValue 7 5 "  -1.95%  "
Set-TextValue 8 5 "  +0.04%  "
Set-TextValue 9 5 "  +0.22%  "
Set-TextValue 10 5 "  -0.99%  "
Set-TextValue 11 4 "6.49"
Set-TextValue 11 5 "  +1.00%  "
Set-TextValue 12 5 "  +0.65%  "
Set-TextValue 13 5 "  +9.08%  "
Set-TextValue 14 4 "36.96"
Set-TextValue 14 5 "  -0.49%  "
Set-TextValue 15 4 "4.449.74"
Set-TextValue 15 5 "  -1.76%  "
Set-TextValue 16 4 "3.823.82"
Set-TextValue 16 5 "  -1.58%  "
Set-TextValue 17 2 "Chainlink"
Set-TextValue 17 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue 17 4 "18.55"
Set-TextValue 17 5 "  +1.69%  "
Set-TextValue 18 2 "WrappedBTC"
Set-TextValue 18 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue 18 4 "68.032.74"
Set-TextValue 18 5 "  -0.40%  "
Set-TextValue 19 4 "7.44"
Set-TextValue 19 5 "  +0.12%  "
Set-TextValue 20 5 "  +0.11%  "
Set-TextValue 21 4 "10.90"
Set-TextValue 21 5 "  +0.51%  "
Set-TextValue 22 4 "470.75"
Set-TextValue 22 5 "  -0.61%  "
Set-TextValue 23 4 "0.738"
Set-TextValue 23 5 "  -0.09%  "
Set-TextValue 24 4 "0.0000152"
Set-TextValue 24 5 "  -8.54%  "
Set-TextValue 25 4 "83.44"
Set-TextValue 25 5 "  -0.48%  "
Set-TextValue 26 4 "2.31"
Set-TextValue 26 5 "  +2.52%  "
Set-TextValue 27 4 "12.19"
Set-TextValue 27 5 "  -0.75%  "
Set-TextValue 28 4 "10.27"
Set-TextValue 28 5 "  +2.67%  "
Set-TextValue 29 5 "  -0.07%  "
Set-TextValue 30 5 "  -1.25%  "
Set-TextValue 31 4 "3.957.10"
Set-TextValue 31 5 "  -1.86%  "
Set-TextValue 32 4 "7.74"
Set-TextValue 32 5 "  -2.09%  "
Set-TextValue 33 5 "  -1.32%  "
Set-TextValue 34 4 "30.78"
Set-TextValue 34 5 "  -2.45%  "
Set-TextValue 35 4 "9.34"
Set-TextValue 35 5 "  -0.83%  "
Set-TextValue 36 4 "3.776.57"
Set-TextValue 36 5 "  -2.10%  "
Set-TextValue 37 5 "  +1.94%  "
Set-TextValue 38 4 "3.81"
Set-TextValue 38 5 "  +5.24%  "
Set-TextValue 39 4 "6.00"
Set-TextValue 39 5 "  +1.08%  "
Set-TextValue 40 5 "  -1.20%  "
Set-TextValue 41 5 "  -2.06%  "
Set-TextValue 42 5 "  -0.09%  "
Set-TextValue 43 4 "0.319"
Set-TextValue 43 5 "  +1.66%  "
Set-TextValue 45 4 "8.81"
Set-TextValue 45 5 "  +1.63%  "
Set-TextValue 46 4 "1.97"
Set-TextValue 46 5 "  -1.32%  "
Set-TextValue 47 4 "410.88"
Set-TextValue 47 5 "  -5.13%  "
Set-TextValue 48 2 "OKB"
Set-TextValue 48 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 48 4 "46.41"
Set-TextValue 48 5 "  -1.95%  "
Set-TextValue 49 2 "FLOKI"
Set-TextValue 49 3 "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue 49 4 "0.000286"
Set-TextValue 49 5 "  -5.45%  "
Set-TextValue 50 2 "VeChain"
Set-TextValue 50 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 50 4 "0.0360"
Set-TextValue 50 5 "  +0.23%  "
Set-TextValue 51 2 "Monero"
Set-TextValue 51 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 51 4 "141.96"
Set-TextValue 51 5 "  -1.12%  "
